# Auto-generated edit script applying cell value updates per the commit diff.
# Each worksheet is selected by name, then specific cells are updated to their new values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 228.33333
$ws.Range("I2").Value = 230
$ws.Range("J2").Value = 225
$ws.Range("K2").Value = 230
$ws.Range("L2").Value = 225
$ws.Range("M2").Value = -117
$ws.Range("N2").Value = -451
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -86
$ws.Range("H58").Value = 2227.5625
$ws.Range("I58").Value = 147.11111
$ws.Range("J58").Value = 4902.4287
$ws.Range("K58").Value = 441.33333
$ws.Range("L58").Value = 14707.2861
$ws.Range("M58").Value = -291.33333
$ws.Range("N58").Value = -15007.2861
$ws.Range("H125").Value = 878572.2
$ws.Range("J125").Value = 1283836.2
$ws.Range("L125").Value = 11554525.8
$ws.Range("N125").Value = -11559445.8
$ws.Range("H132").Value = 5594.579
$ws.Range("I132").Value = 4252.7334
$ws.Range("J132").Value = 10626.5
$ws.Range("K132").Value = 12758.2002
$ws.Range("L132").Value = 31879.5
$ws.Range("M132").Value = -10228.2002
$ws.Range("N132").Value = -36939.5
$ws.Range("H137").Value = 2036.6666
$ws.Range("I137").Value = 1417.3914
$ws.Range("J137").Value = 4071.4285
$ws.Range("K137").Value = 4252.174199999999
$ws.Range("L137").Value = 12214.2855
$ws.Range("M137").Value = -1702.174199999999
$ws.Range("N137").Value = -17314.2855
$ws.Range("H138").Value = 2250.05
$ws.Range("I138").Value = 2211
$ws.Range("J138").Value = 2265.4883
$ws.Range("K138").Value = 6633
$ws.Range("L138").Value = 6796.4649
$ws.Range("M138").Value = -1493
$ws.Range("N138").Value = -17076.4649
$ws.Range("H141").Value = 2387.9333
$ws.Range("I141").Value = 2309.1538
$ws.Range("J141").Value = 2900
$ws.Range("K141").Value = 6927.4614
$ws.Range("L141").Value = 8700
$ws.Range("M141").Value = -1747.4614
$ws.Range("N141").Value = -19060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13665.82
$ws.Range("I32").Value = 5283.8945
$ws.Range("J32").Value = 24776.744
$ws.Range("K32").Value = 5283.8945
$ws.Range("L32").Value = 24776.744
$ws.Range("M32").Value = -4996.8945
$ws.Range("N32").Value = -25350.744
$ws.Range("H74").Value = 1544.4186
$ws.Range("I74").Value = 1038.1177
$ws.Range("J74").Value = 3457.111
$ws.Range("K74").Value = 1038.1177
$ws.Range("L74").Value = 3457.111
$ws.Range("M74").Value = -164.1177
$ws.Range("N74").Value = -5205.111
$ws.Range("H77").Value = 1544.4186
$ws.Range("I77").Value = 1038.1177
$ws.Range("J77").Value = 3457.111
$ws.Range("K77").Value = 5190.5885
$ws.Range("L77").Value = 17285.555
$ws.Range("M77").Value = -822.5884999999998
$ws.Range("N77").Value = -26021.555
$ws.Range("H88").Value = 2166.6667
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 2166.6667
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H95").Value = 350000
$ws.Range("J95").Value = 350000
$ws.Range("L95").Value = 350000
$ws.Range("N95").Value = -355492
$ws.Range("H132").Value = 2312.2144
$ws.Range("I132").Value = 2025.6666
$ws.Range("J132").Value = 4031.5
$ws.Range("K132").Value = 6076.9998
$ws.Range("L132").Value = 12094.5
$ws.Range("M132").Value = -3546.9998
$ws.Range("N132").Value = -17154.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2169.85
$ws.Range("I86").Value = 2217.4707
$ws.Range("J86").Value = 1900
$ws.Range("K86").Value = 2217.4707
$ws.Range("L86").Value = 1900
$ws.Range("M86").Value = -1094.4707
$ws.Range("N86").Value = -4146
$ws.Range("H89").Value = 2169.85
$ws.Range("I89").Value = 2217.4707
$ws.Range("J89").Value = 1900
$ws.Range("K89").Value = 11087.3535
$ws.Range("L89").Value = 9500
$ws.Range("M89").Value = -5471.353499999999
$ws.Range("N89").Value = -20732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2632.0178
$ws.Range("I31").Value = 1855.9706
$ws.Range("J31").Value = 3831.3635
$ws.Range("K31").Value = 1855.9706
$ws.Range("L31").Value = 3831.3635
$ws.Range("M31").Value = -1560.9706
$ws.Range("N31").Value = -4421.363499999999
$ws.Range("H34").Value = 2632.0178
$ws.Range("I34").Value = 1855.9706
$ws.Range("J34").Value = 3831.3635
$ws.Range("K34").Value = 1855.9706
$ws.Range("L34").Value = 3831.3635
$ws.Range("M34").Value = -1653.9706
$ws.Range("N34").Value = -4235.363499999999
$ws.Range("H86").Value = 83336430
$ws.Range("I86").Value = 142860260
$ws.Range("K86").Value = 142860260
$ws.Range("M86").Value = -142859137
$ws.Range("H89").Value = 83336430
$ws.Range("I89").Value = 142860260
$ws.Range("K89").Value = 714301300
$ws.Range("M89").Value = -714295684
$ws.Range("H132").Value = 2484.182
$ws.Range("I132").Value = 2479.111
$ws.Range("J132").Value = 2507
$ws.Range("K132").Value = 7437.333
$ws.Range("L132").Value = 7521
$ws.Range("M132").Value = -4907.333
$ws.Range("N132").Value = -12581
$ws.Range("H135").Value = 53512.5
$ws.Range("J135").Value = 53512.5
$ws.Range("L135").Value = 53512.5
$ws.Range("N135").Value = -63652.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 7775.857
$ws.Range("I80").Value = 3640.5
$ws.Range("J80").Value = 9430
$ws.Range("K80").Value = 10921.5
$ws.Range("L80").Value = 28290
$ws.Range("M80").Value = -9985.5
$ws.Range("N80").Value = -30162
$ws.Range("H83").Value = 7775.857
$ws.Range("I83").Value = 3640.5
$ws.Range("J83").Value = 9430
$ws.Range("K83").Value = 32764.5
$ws.Range("L83").Value = 84870
$ws.Range("M83").Value = -28084.5
$ws.Range("N83").Value = -94230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1986975.8
$ws.Range("I80").Value = 2647.5
$ws.Range("J80").Value = 3574438.5
$ws.Range("K80").Value = 2647.5
$ws.Range("L80").Value = 3574438.5
$ws.Range("M80").Value = -1649.5
$ws.Range("N80").Value = -3576434.5
$ws.Range("H83").Value = 1986975.8
$ws.Range("I83").Value = 2647.5
$ws.Range("J83").Value = 3574438.5
$ws.Range("K83").Value = 13237.5
$ws.Range("L83").Value = 17872192.5
$ws.Range("M83").Value = -8245.5
$ws.Range("N83").Value = -17882176.5
$ws.Range("H122").Value = 2100.7
$ws.Range("I122").Value = 2626.75
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 7880.25
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -5430.25
$ws.Range("N122").Value = -10150
$ws.Range("H124").Value = 56375
$ws.Range("J124").Value = 56375
$ws.Range("L124").Value = 56375
$ws.Range("N124").Value = -66195
$ws.Range("H126").Value = 14691.046
$ws.Range("I126").Value = 3827.875
$ws.Range("J126").Value = 20898.572
$ws.Range("K126").Value = 11483.625
$ws.Range("L126").Value = 62695.716
$ws.Range("M126").Value = -9013.625
$ws.Range("N126").Value = -67635.716
$ws.Range("H132").Value = 3852.9167
$ws.Range("I132").Value = 3279.75
$ws.Range("J132").Value = 4999.25
$ws.Range("K132").Value = 9839.25
$ws.Range("L132").Value = 14997.75
$ws.Range("M132").Value = -7309.25
$ws.Range("N132").Value = -20057.75
$ws.Range("H133").Value = 19176
$ws.Range("J133").Value = 19176
$ws.Range("L133").Value = 19176
$ws.Range("N133").Value = -29296
$ws.Range("H134").Value = 10497.667
$ws.Range("J134").Value = 10497.667
$ws.Range("L134").Value = 31493.001
$ws.Range("N134").Value = -36563.001
$ws.Range("H135").Value = 58496.363
$ws.Range("J135").Value = 58496.363
$ws.Range("L135").Value = 58496.363
$ws.Range("N135").Value = -68636.363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = $null
$ws.Range("H46").Value = 1322.4482
$ws.Range("I46").Value = 1356.2916
$ws.Range("J46").Value = 1160
$ws.Range("K46").Value = 1356.2916
$ws.Range("L46").Value = 1160
$ws.Range("M46").Value = -1168.2916
$ws.Range("N46").Value = -1536
$ws.Range("H68").Value = 1502
$ws.Range("J68").Value = 1453.3334
$ws.Range("L68").Value = 1453.3334
$ws.Range("N68").Value = -2951.3334
$ws.Range("H71").Value = 1502
$ws.Range("J71").Value = 1453.3334
$ws.Range("L71").Value = 7266.666999999999
$ws.Range("N71").Value = -14754.667
$ws.Range("H82").Value = 1492.2222
$ws.Range("I82").Value = 1610
$ws.Range("J82").Value = 1433.3334
$ws.Range("K82").Value = 1610
$ws.Range("L82").Value = 1433.3334
$ws.Range("M82").Value = -1249
$ws.Range("N82").Value = -2155.3334
$ws.Range("H85").Value = 1492.2222
$ws.Range("I85").Value = 1610
$ws.Range("J85").Value = 1433.3334
$ws.Range("K85").Value = 1610
$ws.Range("L85").Value = 1433.3334
$ws.Range("M85").Value = -362
$ws.Range("N85").Value = -3929.3334
$ws.Range("H93").Value = 1525.6086
$ws.Range("I93").Value = 1597.1666
$ws.Range("J93").Value = 1268
$ws.Range("K93").Value = 1597.1666
$ws.Range("L93").Value = 1268
$ws.Range("M93").Value = -349.1666
$ws.Range("N93").Value = -3764
$ws.Range("H103").Value = 100000
$ws.Range("J103").Value = 100000
$ws.Range("L103").Value = 100000
$ws.Range("N103").Value = -102344
$ws.Range("H138").Value = 46406
$ws.Range("J138").Value = 46406
$ws.Range("L138").Value = 46406
$ws.Range("N138").Value = -56686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 17155.666
$ws.Range("J109").Value = 17155.666
$ws.Range("L109").Value = 17155.666
$ws.Range("N109").Value = -19929.666
$ws.Range("H126").Value = 1876.7142
$ws.Range("I126").Value = 1605.4
$ws.Range("J126").Value = 2555
$ws.Range("K126").Value = 4816.200000000001
$ws.Range("L126").Value = 7665
$ws.Range("M126").Value = -2346.200000000001
$ws.Range("N126").Value = -12605
$ws.Range("H136").Value = 3020.7334
$ws.Range("I136").Value = 3534.5557
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 10603.6671
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -8053.667099999999
$ws.Range("N136").Value = -11850
